$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Unified_table")
$ws1.Activate()
$win = $excel.ActiveWindow()
Write-Host "View before:" $win.View()
$win.Zoom = 55
Write-Host "Zoom after:" $win.Zoom()
